{"js": "const replacements = [\n  [\"58\u00d743=2494\", \"36\u00d723=828\"],\n  [\"98\u00d759=5782\", \"25\u00d793=2325\"],\n  [\"85\u00d727=2295\", \"52\u00d775=3900\"],\n  [\"37\u00d720=740\", \"83\u00d741=3403\"],\n  [\"69\u00d711=759\", \"31\u00d736=1116\"],\n  [\"25\u00d739=975\", \"80\u00d720=1600\"],\n  [\"84\u00d726=2184\", \"46\u00d781=3726\"],\n  [\"58\u00d759=3422\", \"18\u00d728=504\"],\n  [\"24\u00d799=2376\", \"18\u00d767=1206\"],\n  [\"86\u00d717=1462\", \"72\u00d773=5256\"],\n  [\"73\u00d790=6570\", \"71\u00d797=6887\"],\n  [\"51\u00d715=765\", \"92\u00d737=3404\"],\n  [\"79\u00d791=7189\", \"29\u00d719=551\"],\n  [\"29\u00d726=754\", \"60\u00d776=4560\"],\n  [\"34\u00d750=1700\", \"48\u00d753=2544\"],\n  [\"24\u00d712=288\", \"75\u00d750=3750\"],\n  [\"62\u00d791=5642\", \"73\u00d737=2701\"],\n  [\"68\u00d759=4012\", \"23\u00d730=690\"],\n  [\"78\u00d779=6162\", \"42\u00d764=2688\"],\n  [\"88\u00d741=3608\", \"65\u00d765=4225\"],\n  [\"98\u00d742=4116\", \"99\u00d789=8811\"],\n  [\"84\u00d739=3276\", \"84\u00d754=4536\"],\n  [\"41\u00d791=3731\", \"88\u00d737=3256\"],\n  [\"63\u00d764=4032\", \"65\u00d734=2210\"],\n  [\"96\u00d748=4608\", \"25\u00d714=350\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\n\nreturn `replaced ${totalReplaced} of ${replacements.length} expected`;", "ps1": "# Replace each two-digit multiplication problem/answer with its updated value.\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{ Old = \"58\u00d743=2494\"; New = \"36\u00d723=828\" },\n  @{ Old = \"98\u00d759=5782\"; New = \"25\u00d793=2325\" },\n  @{ Old = \"85\u00d727=2295\"; New = \"52\u00d775=3900\" },\n  @{ Old = \"37\u00d720=740\"; New = \"83\u00d741=3403\" },\n  @{ Old = \"69\u00d711=759\"; New = \"31\u00d736=1116\" },\n  @{ Old = \"25\u00d739=975\"; New = \"80\u00d720=1600\" },\n  @{ Old = \"84\u00d726=2184\"; New = \"46\u00d781=3726\" },\n  @{ Old = \"58\u00d759=3422\"; New = \"18\u00d728=504\" },\n  @{ Old = \"24\u00d799=2376\"; New = \"18\u00d767=1206\" },\n  @{ Old = \"86\u00d717=1462\"; New = \"72\u00d773=5256\" },\n  @{ Old = \"73\u00d790=6570\"; New = \"71\u00d797=6887\" },\n  @{ Old = \"51\u00d715=765\"; New = \"92\u00d737=3404\" },\n  @{ Old = \"79\u00d791=7189\"; New = \"29\u00d719=551\" },\n  @{ Old = \"29\u00d726=754\"; New = \"60\u00d776=4560\" },\n  @{ Old = \"34\u00d750=1700\"; New = \"48\u00d753=2544\" },\n  @{ Old = \"24\u00d712=288\"; New = \"75\u00d750=3750\" },\n  @{ Old = \"62\u00d791=5642\"; New = \"73\u00d737=2701\" },\n  @{ Old = \"68\u00d759=4012\"; New = \"23\u00d730=690\" },\n  @{ Old = \"78\u00d779=6162\"; New = \"42\u00d764=2688\" },\n  @{ Old = \"88\u00d741=3608\"; New = \"65\u00d765=4225\" },\n  @{ Old = \"98\u00d742=4116\"; New = \"99\u00d789=8811\" },\n  @{ Old = \"84\u00d739=3276\"; New = \"84\u00d754=4536\" },\n  @{ Old = \"41\u00d791=3731\"; New = \"88\u00d737=3256\" },\n  @{ Old = \"63\u00d764=4032\"; New = \"65\u00d734=2210\" },\n  @{ Old = \"96\u00d748=4608\"; New = \"25\u00d714=350\" }\n)\n\nforeach ($pair in $pairs) {\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute(\n    $pair.Old, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false,\n    $pair.New, $wdReplaceAll\n  ) | Out-Null\n}\n"}
